$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("D1").Value = "with zeroingline 437"

# New column D values (row 2 .. row 22)
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0.030425600000000001
$ws.Range("D5").Value = 0.14466300000000001
$ws.Range("D6").Value = 0.26702199999999998
$ws.Range("D7").Value = 0.43737599999999999
$ws.Range("D8").Value = 0.64298900000000003
$ws.Range("D9").Value = 0.87731800000000004
$ws.Range("D10").Value = 1.1409499999999999
$ws.Range("D11").Value = 1.43167
$ws.Range("D12").Value = 1.7473399999999999
$ws.Range("D13").Value = 2.0727799999999998
$ws.Range("D14").Value = 2.4294500000000001
$ws.Range("D15").Value = 2.7867500000000001
$ws.Range("D16").Value = 3.17991
$ws.Range("D17").Value = 3.5922200000000002
$ws.Range("D18").Value = 3.9779399999999998
$ws.Range("D19").Value = 4.3818000000000001
$ws.Range("D20").Value = 4.8424800000000001
$ws.Range("D21").Value = 5.2751700000000001
$ws.Range("D22").Value = 5.7745300000000004

# Update selection to D6 (single cell), matching the diff
$ws.Range("D6").Select()
